# library_J.PLAGGENBERG_07.01.19.xlsx -- final sweep
# 1) index2Name "E7420" -> "E7420L" (column K, rows 2-49)
# 2) column L (roboticLibraryPrep) literal FALSE -> formula =FALSE() (rows 2-49)
# 3) move the active selection from column L to column K, scrolled down a bit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 11).Value   = "E7420L"
    $ws.Cells.Item($row, 12).Formula = "=FALSE()"
}

# Scroll the window so row 13 / column F sit in the top-left corner, then
# restore the selection onto the (now current) data column K.
$win = $excel.ActiveWindow
$win.ScrollRow    = 13
$win.ScrollColumn = 6

$ws.Range("K2:K49").Select()
